$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated Team Belbin's Analysis skill levels (column C, "Skill Level (1-5)")
$ws.Range("C7").Value  = 3
$ws.Range("C9").Value  = 3.8
$ws.Range("C11").Value = 4.6
$ws.Range("C12").Value = 4
$ws.Range("C13").Value = 4
$ws.Range("C15").Value = 2.8
$ws.Range("C16").Value = 3.6
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 2.8
$ws.Range("C19").Value = 2.8
$ws.Range("C20").Value = 3.4
$ws.Range("C21").Value = 4.4
$ws.Range("C22").Value = 3.8
$ws.Range("C23").Value = 3

# Restore the view/selection state recorded when the file was last saved
$ws.Range("A11").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()
